$wb = $excel.ActiveWorkbook

# Logs sheet - add new row 14 (copy of the "Afmelding nieuwsbrief" entry)
$logs = $wb.Worksheets.Item("Logs")
$logs.Range("A14").Value = "Afmelding nieuwsbrief"
$logs.Range("B14").Value = "mailmind.test@zohomail.eu"
$logs.Range("C14").Value = "Graag afmelden voor de nieuwsbrief. Dank u."
$logs.Range("D14").Value = "Afmelding / Nieuwsbrief"
$logs.Range("F14").Value = "2025-06-20 14:00:11"
$logs.Range("G14").Value = "Nee"

# Dashboard sheet - increment count for "Afmelding / Nieuwsbrief"
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B3").Value = 3

# Extend conditional formatting ranges to include the new row 14
$dFcs = $logs.Range("D2:D13").FormatConditions
for ($i = 1; $i -le $dFcs.Count; $i++) {
    $dFcs.Item($i).ModifyAppliesToRange($logs.Range("D2:D14"))
}

$gFcs = $logs.Range("G2:G13").FormatConditions
for ($i = 1; $i -le $gFcs.Count; $i++) {
    $gFcs.Item($i).ModifyAppliesToRange($logs.Range("G2:G14"))
}
